$wb = $excel.ActiveWorkbook

# --- Veda sheet: column C scenario-category helper column ---
# C8 used to hold a hard-coded shared-string ("C2"); it is retyped as a
# formula that copies the cell above (=C7), and that formula is then
# filled down through C17, turning the whole C8:C17 block into a
# relative "=previous row" chain that is ultimately anchored on C7
# ("C1"). This changes the effective lookup value used throughout the
# sheet (Veda!$C$5, via VLOOKUP(B3,...)) from "C3" to "C1", which in
# turn recalculates every SUMIFS() pulling from ar6_r10 that filters on
# Veda!$C$5.
$wsVeda = $wb.Worksheets.Item("Veda")
$wsVeda.Range("C8").Formula = "=C7"
for ($r = 9; $r -le 17; $r++) {
    $prevRow = $r - 1
    $wsVeda.Range("C$r").Formula = "=C$prevRow"
}

# --- View-state bookkeeping: which sheet/tab is active & what's selected ---
$wsFuel = $wb.Worksheets.Item("fuel_prices")
$wsIea  = $wb.Worksheets.Item("iea_data")
$wsAr6  = $wb.Worksheets.Item("ar6_r10")
$wsEv   = $wb.Worksheets.Item("ev_charging_uc")

$wsFuel.Range("A1:H1").Select()
$wsIea.Range("A1:H1").Select()
$wsAr6.Range("A1:K1").Select()

# ev_charging_uc becomes the active/visible tab (activeTab moves from
# index 3 -> 4), keeping its existing H21 selection.
$wsEv.Activate()
$wsEv.Range("H21").Select()
